{"js": "// Replace each two-digit-by-two-digit multiplication prompt in the table\n// with its new pair of operands. Every \"old\" string is unique in the\n// document and every \"new\" string is distinct from all \"old\" strings, so\n// the replacements can be applied independently/sequentially in any order.\nconst replacements = [\n  [\"12\u00d791=\", \"25\u00d747=\"],\n  [\"47\u00d792=\", \"41\u00d747=\"],\n  [\"69\u00d797=\", \"23\u00d738=\"],\n  [\"53\u00d720=\", \"86\u00d796=\"],\n  [\"31\u00d727=\", \"71\u00d720=\"],\n  [\"45\u00d721=\", \"65\u00d788=\"],\n  [\"80\u00d783=\", \"89\u00d782=\"],\n  [\"51\u00d759=\", \"79\u00d774=\"],\n  [\"98\u00d760=\", \"59\u00d766=\"],\n  [\"93\u00d796=\", \"31\u00d749=\"],\n  [\"64\u00d790=\", \"44\u00d746=\"],\n  [\"38\u00d762=\", \"36\u00d758=\"],\n  [\"11\u00d732=\", \"87\u00d793=\"],\n  [\"78\u00d791=\", \"56\u00d768=\"],\n  [\"71\u00d778=\", \"68\u00d712=\"],\n  [\"11\u00d713=\", \"31\u00d777=\"],\n  [\"34\u00d746=\", \"31\u00d763=\"],\n  [\"57\u00d712=\", \"90\u00d735=\"],\n  [\"67\u00d733=\", \"28\u00d727=\"],\n  [\"53\u00d779=\", \"20\u00d774=\"],\n  [\"64\u00d713=\", \"66\u00d788=\"],\n  [\"94\u00d746=\", \"52\u00d752=\"],\n  [\"96\u00d797=\", \"44\u00d717=\"],\n  [\"23\u00d791=\", \"88\u00d732=\"],\n  [\"24\u00d795=\", \"54\u00d780=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const found = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  await context.sync();\n\n  for (const range of found.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace each two-digit-by-two-digit multiplication prompt in the table\n# with its new pair of operands. Every \"old\" string is unique in the\n# document and every \"new\" string is distinct from all \"old\" strings, so\n# the replacements can be applied independently/sequentially in any order.\n$d = $word.ActiveDocument\n\n$replacements = @(\n  @(\"12\u00d791=\", \"25\u00d747=\"),\n  @(\"47\u00d792=\", \"41\u00d747=\"),\n  @(\"69\u00d797=\", \"23\u00d738=\"),\n  @(\"53\u00d720=\", \"86\u00d796=\"),\n  @(\"31\u00d727=\", \"71\u00d720=\"),\n  @(\"45\u00d721=\", \"65\u00d788=\"),\n  @(\"80\u00d783=\", \"89\u00d782=\"),\n  @(\"51\u00d759=\", \"79\u00d774=\"),\n  @(\"98\u00d760=\", \"59\u00d766=\"),\n  @(\"93\u00d796=\", \"31\u00d749=\"),\n  @(\"64\u00d790=\", \"44\u00d746=\"),\n  @(\"38\u00d762=\", \"36\u00d758=\"),\n  @(\"11\u00d732=\", \"87\u00d793=\"),\n  @(\"78\u00d791=\", \"56\u00d768=\"),\n  @(\"71\u00d778=\", \"68\u00d712=\"),\n  @(\"11\u00d713=\", \"31\u00d777=\"),\n  @(\"34\u00d746=\", \"31\u00d763=\"),\n  @(\"57\u00d712=\", \"90\u00d735=\"),\n  @(\"67\u00d733=\", \"28\u00d727=\"),\n  @(\"53\u00d779=\", \"20\u00d774=\"),\n  @(\"64\u00d713=\", \"66\u00d788=\"),\n  @(\"94\u00d746=\", \"52\u00d752=\"),\n  @(\"96\u00d797=\", \"44\u00d717=\"),\n  @(\"23\u00d791=\", \"88\u00d732=\"),\n  @(\"24\u00d795=\", \"54\u00d780=\")\n)\n\nforeach ($pair in $replacements) {\n  $oldText = $pair[0]\n  $newText = $pair[1]\n\n  $find = $d.Content.Find\n  $find.Text = $oldText\n  $find.Replacement.Text = $newText\n  $find.Execute($null, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null\n}\n"}
